$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("D1")
$scratch.HorizontalAlignment = -4131
$scratch.VerticalAlignment = -4108

$ws.Range("A8").ClearContents()
$ws.Range("A7:A8").Merge($false)

$scratch.Copy()
$ws.Range("A7:A8").PasteSpecial(-4122)

$ws.Range("D1").Clear()
